# Updated Agile Document Excel sheets with latest data
# - "Assigned To" column (F) on the Defects sheet is reassigned to a
#   single person ("Aryan Khurana") for rows 3-10, replacing the various
#   names that were there before.
# - The active selection on the Defects sheet moves to F10.
# - Column F is widened slightly to fit the new, longer name.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Defects")

# Reassign every defect in the "Assigned To" column to Aryan Khurana.
$ws.Range("F3:F10").Value = "Aryan Khurana"

# Widen column F so the new name fits (closest attainable width).
$ws.Columns("F").ColumnWidth = 14

# Move / leave the active cell selection on F10.
[void]$ws.Range("F10").Select()
